$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.194.40'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +3.41%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.913.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +2.78%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -1.49%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''314.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.64%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -1.55%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4865'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +1.41%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3840'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +3.04%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07399'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.99%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.9515'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +1.55%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''21.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +3.88%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.07805'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.84%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.915.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +2.04%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''5.568'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +2.63%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''6.649'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.81%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''92.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +2.40%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''1.006'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -1.53%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.000008909'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.84%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''1.005'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -1.38%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''28.190.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +3.26%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''15.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +1.62%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''5.164'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +1.01%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''2.140.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +1.63%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  +2.71%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''156.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +1.96%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''1.935'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.77%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''18.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.86%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''2.098'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +4.85%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''116.34'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +0.53%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''5.019'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.59%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.08900'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.06%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.319'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -0.87%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''1.252'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +5.08%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.7806'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +5.47%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''4.691'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +2.56%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''2.767'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +3.29%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = '''VeChain'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = '''0.02055'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +1.04%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = '''TrustWalletToken'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = '''1.127'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +0.29%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.5632'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +5.05%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.05380'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +2.28%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''3.029'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.02%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''7.090'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -0.52%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''8.620'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +3.25%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.1535'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +0.09%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = '''EnergySwap'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''10.88'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +2.60%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''Decentraland'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''0.4945'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +2.97%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''106.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +3.19%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''NEARProtocol'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''1.684'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +2.95%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = '''PaxDollar'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''1.003'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -1.68%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''69.23'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +3.95%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.06145'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +1.12%  '
$ws.Range('E51').Style = 'Normal'
